$wb = $excel.ActiveWorkbook

# ALC!row28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1130.1538
$ws.Range("I28").Value = 1273.091
$ws.Range("J28").Value = 344
$ws.Range("K28").Value = 1273.091
$ws.Range("L28").Value = 344
$ws.Range("M28").Value = -788.0909999999999
$ws.Range("N28").Value = -1314

# ALC!row32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1362.4
$ws.Range("I32").Value = 998.6667
$ws.Range("J32").Value = 1518.2858
$ws.Range("K32").Value = 998.6667
$ws.Range("L32").Value = 1518.2858
$ws.Range("M32").Value = -672.6667
$ws.Range("N32").Value = -2170.2858

# ALC!row43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2137
$ws.Range("I43").Value = 2158.3333
$ws.Range("J43").Value = 2121
$ws.Range("K43").Value = 2158.3333
$ws.Range("L43").Value = 2121
$ws.Range("M43").Value = -2089.3333
$ws.Range("N43").Value = -2259

# ALC!row64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5386.5
$ws.Range("I64").Value = 5495.75
$ws.Range("J64").Value = 4949.5
$ws.Range("K64").Value = 5495.75
$ws.Range("L64").Value = 4949.5
$ws.Range("M64").Value = -5247.75
$ws.Range("N64").Value = -5445.5

# ALC!row67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 5386.5
$ws.Range("I67").Value = 5495.75
$ws.Range("J67").Value = 4949.5
$ws.Range("K67").Value = 5495.75
$ws.Range("L67").Value = 4949.5
$ws.Range("M67").Value = -4637.75
$ws.Range("N67").Value = -6665.5

# ALC!row88
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 2249
$ws.Range("J88").Value = 1999.5
$ws.Range("L88").Value = 1999.5
$ws.Range("N88").Value = -2811.5

# ALC!row91
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 2249
$ws.Range("J91").Value = 1999.5
$ws.Range("L91").Value = 1999.5
$ws.Range("N91").Value = -4807.5

# ALC!row98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2261
$ws.Range("I98").Value = 1770.2307
$ws.Range("J98").Value = 4387.6665
$ws.Range("K98").Value = 1770.2307
$ws.Range("L98").Value = 4387.6665
$ws.Range("M98").Value = -272.2307000000001
$ws.Range("N98").Value = -7383.6665

# ALC!row111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 2089.818
$ws.Range("I111").Value = 1617.1818
$ws.Range("K111").Value = 4851.5454
$ws.Range("M111").Value = -1784.5454

# ALC!row116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 16925.924
$ws.Range("I116").Value = 25473.75
$ws.Range("J116").Value = 13126.889
$ws.Range("K116").Value = 25473.75
$ws.Range("L116").Value = 13126.889
$ws.Range("M116").Value = -22031.75
$ws.Range("N116").Value = -20010.889

# ALC!row122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2261
$ws.Range("I122").Value = 1770.2307
$ws.Range("J122").Value = 4387.6665
$ws.Range("K122").Value = 5310.6921
$ws.Range("L122").Value = 13162.9995
$ws.Range("M122").Value = -2860.6921
$ws.Range("N122").Value = -18062.9995

# ALC!row131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 3320.818
$ws.Range("I131").Value = 1544.6471
$ws.Range("J131").Value = 9359.799999999999
$ws.Range("K131").Value = 4633.9413
$ws.Range("L131").Value = 28079.4
$ws.Range("M131").Value = 406.0587000000005
$ws.Range("N131").Value = -38159.39999999999

# ALC!row133
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -60120

# ARM!row74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2966.5
$ws.Range("I74").Value = 1878.6111
$ws.Range("J74").Value = 4054.389
$ws.Range("K74").Value = 1878.6111
$ws.Range("L74").Value = 4054.389
$ws.Range("M74").Value = -1004.6111
$ws.Range("N74").Value = -5802.389

# ARM!row77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2966.5
$ws.Range("I77").Value = 1878.6111
$ws.Range("J77").Value = 4054.389
$ws.Range("K77").Value = 9393.0555
$ws.Range("L77").Value = 20271.945
$ws.Range("M77").Value = -5025.0555
$ws.Range("N77").Value = -29007.945

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3152.9546
$ws.Range("I31").Value = 2106.6155
$ws.Range("K31").Value = 2106.6155
$ws.Range("M31").Value = -1811.6155

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3152.9546
$ws.Range("I34").Value = 2106.6155
$ws.Range("K34").Value = 2106.6155
$ws.Range("M34").Value = -1904.6155

# CRP!row58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5207
$ws.Range("J58").Value = 7347.3335
$ws.Range("L58").Value = 7347.3335
$ws.Range("N58").Value = -7753.3335

# CRP!row99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2725.5557
$ws.Range("I99").Value = 2353.2307
$ws.Range("K99").Value = 2353.2307
$ws.Range("M99").Value = -855.2307000000001

# CRP!row126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2725.5557
$ws.Range("I126").Value = 2353.2307
$ws.Range("K126").Value = 7059.6921
$ws.Range("M126").Value = -4589.6921

# CRP!row136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 5207
$ws.Range("J136").Value = 7347.3335
$ws.Range("L136").Value = 22042.0005
$ws.Range("N136").Value = -27142.0005

# CUL!row75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 4998.5
$ws.Range("I75").Value = 4999
$ws.Range("J75").Value = 4998
$ws.Range("K75").Value = 14997
$ws.Range("L75").Value = 14994
$ws.Range("M75").Value = -13999
$ws.Range("N75").Value = -16990

# CUL!row78
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 4998.5
$ws.Range("I78").Value = 4999
$ws.Range("J78").Value = 4998
$ws.Range("K78").Value = 44991
$ws.Range("L78").Value = 44982
$ws.Range("M78").Value = -39999
$ws.Range("N78").Value = -54966

# CUL!row131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4997397.5
$ws.Range("J131").Value = 4119116
$ws.Range("L131").Value = 12357348
$ws.Range("N131").Value = -12367428

# GSM!row82
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H82").Value = 50000
$ws.Range("I82").Value = 50000
$ws.Range("K82").Value = 50000
$ws.Range("M82").Value = -49617

# GSM!row85
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H85").Value = 50000
$ws.Range("I85").Value = 50000
$ws.Range("K85").Value = 50000
$ws.Range("M85").Value = -48674

# GSM!row113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2262.4285
$ws.Range("I113").Value = 2262.4285
$ws.Range("K113").Value = 2262.4285
$ws.Range("M113").Value = -92.42849999999999

# LTW!row20
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 4169166.8
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 4169166.8
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 4169166.8
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -4169618.8

# WVR!row81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5336.143
$ws.Range("I81").Value = 5336.143
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 10672.286
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -9611.286
$ws.Range("N81").ClearContents()

# WVR!row84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 5336.143
$ws.Range("I84").Value = 5336.143
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 53361.43
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -48057.43
$ws.Range("N84").ClearContents()

# WVR!row126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1852.5
$ws.Range("I126").Value = 1960
$ws.Range("K126").Value = 5880
$ws.Range("M126").Value = -3410
